# "Preço Adicionado na Resposta" — remove the "Preço (R$)" price column
# (header + values) from the medication table. Clearing the contents (while
# leaving the cell formatting/style intact) removes the now-unused
# "Preço (R$)" entry from the shared-string table on save.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F1:F53").ClearContents()
